# Update "上海-漫展信息.xlsx": refresh "想去人数" (want-to-go count) figures
# across all four sheets, and append two newly-scraped 燃梦BACG PRO events
# to the "展览" (Exhibition) sheet (rows 44-45).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition) — update column F ("想去人数") values
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,6).Value = 1405
$ws.Cells.Item(3,6).Value = 101
$ws.Cells.Item(5,6).Value = 6216
$ws.Cells.Item(6,6).Value = 498
$ws.Cells.Item(7,6).Value = 1051
$ws.Cells.Item(8,6).Value = 18
$ws.Cells.Item(9,6).Value = 3466
$ws.Cells.Item(10,6).Value = 6692
$ws.Cells.Item(12,6).Value = 1335
$ws.Cells.Item(13,6).Value = 771
$ws.Cells.Item(14,6).Value = 101
$ws.Cells.Item(17,6).Value = 1120
$ws.Cells.Item(19,6).Value = 113
$ws.Cells.Item(21,6).Value = 177
$ws.Cells.Item(23,6).Value = 996
$ws.Cells.Item(24,6).Value = 323
$ws.Cells.Item(26,6).Value = 19
$ws.Cells.Item(27,6).Value = 110
$ws.Cells.Item(30,6).Value = 20
$ws.Cells.Item(31,6).Value = 64
$ws.Cells.Item(32,6).Value = 4
$ws.Cells.Item(33,6).Value = 21
$ws.Cells.Item(34,6).Value = 21
$ws.Cells.Item(35,6).Value = 4
$ws.Cells.Item(36,6).Value = 316
$ws.Cells.Item(37,6).Value = 23
$ws.Cells.Item(39,6).Value = 295
$ws.Cells.Item(41,6).Value = 2
$ws.Cells.Item(43,6).Value = 97

# Append two new rows (44, 45) at the bottom of the 展览 sheet.
# Copy the formatting of row 43's numbering cell (bold/centered/bordered
# style) onto the new numbering cells A44/A45 before writing their values.
$ws.Range("A43").Copy()
$ws.Range("A44").PasteSpecial(-4122)
$ws.Range("A43").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(44,1).Value = 43
$ws.Cells.Item(44,2).NumberFormat = "@"
$ws.Cells.Item(44,2).Value = "2024-08-17"
$ws.Cells.Item(44,3).Value = "上海·第六届燃梦BACG PRO动漫嘉年华-我们在燃梦相遇吧！"
$ws.Cells.Item(44,4).Value = "盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)"
$ws.Cells.Item(44,5).Value = "2024.08.17 11:00-08.18 16:00"
$ws.Cells.Item(44,6).Value = 2
$ws.Cells.Item(44,7).Value = 65.8
$ws.Cells.Item(44,8).Value = "https://show.bilibili.com/platform/detail.html?id=85239"
$ws.Cells.Item(44,9).Value = "//i1.hdslb.com/bfs/openplatform/202405/mzD4rhY21715109458100.jpeg"

$ws.Cells.Item(45,1).Value = 44
$ws.Cells.Item(45,2).NumberFormat = "@"
$ws.Cells.Item(45,2).Value = "2024-09-15"
$ws.Cells.Item(45,3).Value = "上海·第七届燃梦BACG PRO动漫嘉年华·我们在燃梦相遇吧！"
$ws.Cells.Item(45,4).Value = "漕宝路3366号 七宝万科广场"
$ws.Cells.Item(45,5).Value = "2024.09.15 10:00-09.16 16:00"
$ws.Cells.Item(45,6).Value = 1
$ws.Cells.Item(45,7).Value = 58
$ws.Cells.Item(45,8).Value = "https://show.bilibili.com/platform/detail.html?id=85240"
$ws.Cells.Item(45,9).Value = "//i2.hdslb.com/bfs/openplatform/202405/RhRh4LZM1715110314459.jpeg"

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performance) — update column F ("想去人数") values
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(6,6).Value = 515
$ws.Cells.Item(12,6).Value = 12
$ws.Cells.Item(21,6).Value = 184
$ws.Cells.Item(30,6).Value = 691

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local Life) — update column F ("想去人数") values
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4,6).Value = 702
$ws.Cells.Item(5,6).Value = 828
$ws.Cells.Item(6,6).Value = 571
$ws.Cells.Item(7,6).Value = 283
$ws.Cells.Item(8,6).Value = 1108

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All Types) — update column F ("想去人数") values
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2,6).Value = 702
$ws.Cells.Item(3,6).Value = 828
$ws.Cells.Item(5,6).Value = 101
$ws.Cells.Item(6,6).Value = 571
$ws.Cells.Item(7,6).Value = 283
$ws.Cells.Item(8,6).Value = 515
$ws.Cells.Item(10,6).Value = 6216
$ws.Cells.Item(11,6).Value = 498
$ws.Cells.Item(12,6).Value = 1051
$ws.Cells.Item(13,6).Value = 18
$ws.Cells.Item(15,6).Value = 6692
$ws.Cells.Item(18,6).Value = 1335
$ws.Cells.Item(19,6).Value = 12
$ws.Cells.Item(24,6).Value = 1108
$ws.Cells.Item(27,6).Value = 184
$ws.Cells.Item(29,6).Value = 113
$ws.Cells.Item(33,6).Value = 19
$ws.Cells.Item(36,6).Value = 20
$ws.Cells.Item(37,6).Value = 64
$ws.Cells.Item(38,6).Value = 4
$ws.Cells.Item(41,6).Value = 21
$ws.Cells.Item(43,6).Value = 316
$ws.Cells.Item(46,6).Value = 295
$ws.Cells.Item(49,6).Value = 97

Write-Output "edit complete"
